# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-10-12 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-10-13 Friday", 2)

# Update the division problems in the table, cell-by-cell (table has
# data rows at 1, 5, 9, 13, 17 -- each with 5 columns). Using explicit
# row/column addressing avoids any ambiguity from duplicate problem
# text (e.g. "92÷3=" appears three times with three different
# replacements).
$tbl = $d.Tables.Item(1)

$newValues = @{
    1  = @("99÷5=", "79÷7=", "85÷2=", "97÷8=", "58÷5=")
    5  = @("91÷6=", "79÷3=", "70÷5=", "39÷9=", "85÷8=")
    9  = @("83÷7=", "42÷2=", "69÷3=", "82÷7=", "13÷6=")
    13 = @("78÷5=", "85÷2=", "63÷6=", "39÷6=", "28÷5=")
    17 = @("67÷4=", "34÷2=", "33÷2=", "60÷4=", "70÷3=")
}

foreach ($rowIndex in $newValues.Keys) {
    $values = $newValues[$rowIndex]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $tbl.Cell($rowIndex, $col)
        $cellRange = $cell.Range
        $cellRange.MoveEnd(1, -1) | Out-Null
        $cellRange.Text = $values[$col - 1]
    }
}
